# Generate Report for Handback
# Update the timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-11-29 04:54:17"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn.Range("H2").Value = "2016-11-29 04:54:03"
$wsZhCn.Range("K2").Value = "2016-11-29 04:54:53"

# de-de sheet: "Correspond Handoff Datetime" (mirrors the Overview generate date)
# and "Correspond Handback DateTime"
$wsDeDe.Range("H2").Value = "2016-11-29 04:54:17"
$wsDeDe.Range("K2").Value = "2016-11-29 04:55:13"
